$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C29 with the expanded Work Microwave announcement text (IBC2016)
$ws.Range("C29").Value = 'Work Microwave plans to implement DVB-S2X for both their Broadcast Modulator SDM2, and professional IP Modem, SK-IP. The company will also introduce DVB-S2X products at IBC in September. These two device families make use of the broadcast profile and the professional profile.    HOLZKIRCHEN, Germany — Sept. 20, 2016 — WORK Microwave, a leading European manufacturer of advanced satellite communications equipment, today announced a new single-carrier Wideband Demodulator that, when combined with the company’s Broadcast Modulator, provides satellite operators with one of the first end-to-end wideband transmission and reception solutions. Single-carrier operation mode for the Wideband Broadcast Modulator and Demodulator is based on the DVB-S2X standard, allowing the most efficient statistical multiplexing of different services and transmission over high-throughput satellite transponders with guaranteed power efficiency. Ideal for next-generation, high-speed IP-based broadcast and broadband access applications in Ka-, Ku-, Q-, and V-band satellite systems, WORK Microwave’s end-to-end wideband solution enables symbol rates of up to 500Msps.   WORK Microwave displayed a prototype of the Wideband Demodulator at IBC2016, with commercial availability in Q2 2017. The DVB-S2X Wideband Broadcast Modulator is now shipping.'

# Add the new URL cell for Work Microwave in B29, matching style of neighboring cells
$ws.Range("B29").Value = 'https://work-microwave.com'

# Row 29 grows taller to fit the added text
$ws.Rows.Item(29).RowHeight = 330

# Update the sheet view to match the new scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 40
$ws.Range("B29").Select()
